$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy style from existing header cell (A1) to the new header cells
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Boolean outlier flag data for rows 2-13 (TRUE/FALSE), all FALSE except row 4 which is TRUE
$values = @(
    $false,
    $false,
    $true,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false,
    $false
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
    $ws.Cells.Item($row, 7).Value = $values[$i]
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
